# Update "Forecast Comparison" sheet with corrected forecast output:
# - Insert a new "Week_Start_Date" column after "Week" (new column B)
# - Shift existing columns ASIN..is_holiday_week right by one
# - Update "Week" values from zero-padded (W01..W09) to unpadded (W1..W9)
# - Populate Week_Start_Date with the corresponding week's start date
# - Change ASIN column from a computed/number to inline string (already string) - unaffected
# - Update the is_holiday_week column type to boolean

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Forecast Comparison")

# Insert a new column before column B (ASIN) to make room for Week_Start_Date
$ws.Columns.Item(2).Insert()

# Set header for the new column
$ws.Cells.Item(1, 2).Value = "Week_Start_Date"

# Force the new column to be stored as text so date-like strings are not
# auto-converted into date serial numbers
$ws.Columns.Item(2).NumberFormat = "@"

# Week labels (unpadded) and their corresponding start dates
$weeks = @("W1","W2","W3","W4","W5","W6","W7","W8","W9","W10","W11","W12","W13","W14","W15","W16")
$startDates = @(
    "2025-01-05","2025-01-12","2025-01-19","2025-01-26",
    "2025-02-02","2025-02-09","2025-02-16","2025-02-23",
    "2025-03-02","2025-03-09","2025-03-16","2025-03-23",
    "2025-03-30","2025-04-06","2025-04-13","2025-04-20"
)

for ($i = 0; $i -lt $weeks.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 1).Value = $weeks[$i]
    $ws.Cells.Item($row, 2).Value = $startDates[$i]
    # is_holiday_week is now a boolean column (was numeric) - write boolean FALSE
    $ws.Cells.Item($row, 10).Value = $false
}

# Update dimension-related bookkeeping is handled automatically by the runtime.
